# Adds a "2022-Q3" quarter: a new worksheet with the per-fund holdings,
# plus a new summary row on the "总计" sheet.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$val) {
    # Force literal text storage (avoids Excel re-parsing numeric-looking
    # strings like "001645" or "36.92" into numbers and losing the
    # leading zeros / exact text form).
    $range.Value = "'" + $val
}

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet, positioned right before "2022-Q2"
#    (this naturally shifts 2022-Q2 .. 2020-Q4 one slot later).
# ---------------------------------------------------------------------
$anchor = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($anchor)
$newSheet.Name = "2022-Q3"

# Re-fetch worksheet references now that the sheet collection changed -
# handles pasted-from COM object (captured before Add()) silently lose
# their formatting-copy ability once a new sheet is inserted.
$anchor = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Item("2022-Q3")

# Copy header-row look & feel (bold / border / centered) from the
# existing "2022-Q2" header row.
$anchor.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Fund-level holdings for 2022-Q3.
# code, name, scale, stockPosTotal, posPct, heldValue, rank
$funds = @(
    @("001645", "国泰大健康股票A", "36.92", "92.54", "7.48", "2.7616", 5),
    @("020001", "国泰金鹰增长灵活配置混合", "20.78", "94.53", "8.20", "1.7040", 7),
    @("009804", "国泰研究优势混合", "12.23", "92.87", "6.99", "0.8549", 8),
    @("009805", "国泰医药健康股票A", "8.71", "94.26", "8.50", "0.7404", 4),
    @("160215", "国泰价值经典灵活配置混合（LOF）", "6.01", "93.77", "7.93", "0.4766", 5),
    @("011321", "国泰大健康股票C", "3.20", "92.54", "7.48", "0.2394", 5),
    @("011738", "华安兴安优选一年持有期混合A", "13.36", "55.44", "1.67", "0.2231", 9),
    @("008370", "国泰研究精选两年持有期混合", "2.98", "93.29", "7.47", "0.2226", 7),
    @("011739", "华安兴安优选一年持有期混合C", "8.27", "55.44", "1.67", "0.1381", 9),
    @("005695", "华安睿明两年定期开放灵活配置混合A", "4.27", "93.55", "2.96", "0.1264", 9),
    @("005585", "银河文体娱乐主题灵活配置混合A", "3.01", "90.28", "4.04", "0.1216", 9),
    @("014786", "惠升品质优选混合A", "1.51", "82.37", "6.53", "0.0986", 2),
    @("015110", "惠升领先优选混合A", "1.50", "81.38", "6.47", "0.0970", 2),
    @("011326", "国泰医药健康股票C", "0.52", "94.26", "8.50", "0.0442", 4),
    @("015667", "银河文体娱乐主题灵活配置混合C", "0.41", "90.28", "4.04", "0.0166", 9),
    @("009409", "华安添福18个月持有期混合A", "0.72", "21.56", "0.94", "0.0068", 9),
    @("005696", "华安睿明两年定期开放灵活配置混合C", "0.07", "93.55", "2.96", "0.0021", 9),
    @("009410", "华安添福18个月持有期混合C", "0.09", "21.56", "0.94", "0.0008", 9),
    @("015111", "惠升领先优选混合C", "0.00", "81.38", "6.47", 0, 2),
    @("014787", "惠升品质优选混合C", "0.00", "82.37", "6.53", 0, 2)
)

for ($i = 0; $i -lt $funds.Length; $i++) {
    $row = $i + 2
    $f = $funds[$i]

    # Row index style (bold / border / centered) matches the other
    # quarter sheets' column A.
    $anchor.Cells.Item($i + 2, 1).Copy()
    $newSheet.Cells.Item($row, 1).PasteSpecial(-4122)
    $excel.CutCopyMode = $false
    $newSheet.Cells.Item($row, 1).Value = $i

    Set-TextValue $newSheet.Cells.Item($row, 2) $f[0]
    $newSheet.Cells.Item($row, 3).Value = $f[1]
    Set-TextValue $newSheet.Cells.Item($row, 4) $f[2]
    Set-TextValue $newSheet.Cells.Item($row, 5) $f[3]
    Set-TextValue $newSheet.Cells.Item($row, 6) $f[4]

    if ($f[5] -is [string]) {
        Set-TextValue $newSheet.Cells.Item($row, 7) $f[5]
    } else {
        $newSheet.Cells.Item($row, 7).Value = $f[5]
    }

    $newSheet.Cells.Item($row, 8).Value = $f[6]
}

# ---------------------------------------------------------------------
# 2) Add the new summary row into "总计". The quarter rows (B/C/D) move
#    down by one so the newest quarter lands on top at row 2, but the
#    leading index column (A) is untouched for existing rows and simply
#    grows by one more sequential value for the new last row - matching
#    the source data's append-only index semantics.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Snapshot existing quarter rows (2..8) top-to-bottom before overwriting.
$oldQuarters = @()
for ($r = 2; $r -le 8; $r++) {
    $oldQuarters += , @(
        $totalSheet.Cells.Item($r, 2).Value(),
        $totalSheet.Cells.Item($r, 3).Value(),
        $totalSheet.Cells.Item($r, 4).Value()
    )
}

# Shift them down one row: old row r -> new row r+1.
for ($i = $oldQuarters.Length - 1; $i -ge 0; $i--) {
    $destRow = $i + 3
    $vals = $oldQuarters[$i]
    $totalSheet.Cells.Item($destRow, 2).Value = $vals[0]
    $totalSheet.Cells.Item($destRow, 3).Value = $vals[1]
    $totalSheet.Cells.Item($destRow, 4).Value = $vals[2]
}

# New index cell for the appended row (row 9), matching column A's style.
$totalSheet.Cells.Item(8, 1).Copy()
$totalSheet.Cells.Item(9, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$totalSheet.Cells.Item(9, 1).Value = 7

# New quarter's data goes on top, row 2 (A2 stays 0, untouched).
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 20
$totalSheet.Cells.Item(2, 4).Value = 7.87
